$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132076025009155
$ws.Range("B1").Value = 2.235070705413818
$ws.Range("C1").Value = 10.43289947509766
$ws.Range("D1").Value = 2.238921642303467
$ws.Range("E1").Value = 1.282428503036499
